$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values on the rows that were repulled/recalculated
$ws.Range("F4").Value = -2
$ws.Range("F11").Value = 3
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = 0
$ws.Range("F23").Value = -5
$ws.Range("F27").Value = -6
$ws.Range("F31").Value = -1
$ws.Range("F33").Value = -8
